# Apply the changes described by the commit:
# - Rename sheets "wt" -> "wt_log2_expression" and "dcin5" -> "dcin5_log2_expression"
# - Make "dcin5_log2_expression" the active/selected sheet (was "optimization_parameters")

$wb = $excel.ActiveWorkbook

# Rename the worksheets
$wsWt = $wb.Worksheets.Item("wt")
$wsWt.Name = "wt_log2_expression"

$wsDcin5 = $wb.Worksheets.Item("dcin5")
$wsDcin5.Name = "dcin5_log2_expression"

# Make the dcin5_log2_expression sheet the active tab (clears tabSelected on the
# previously active sheet and sets it on this one; also updates workbook's
# activeTab / removes firstSheet scrolling offset)
$wsDcin5.Activate()
$wsDcin5.Select()
